# Reorders the Notation (column B) rows within each Collection group on the
# "Collections" sheet, carrying each row's PrefLabel/Definition/UsageNote/
# ScopeNote (columns C-F) along with it, per the target permutation.
#
# Mapping is expressed as: destination row -> source row (both 2..50),
# i.e. after the edit, row <dest> contains what row <source> contained
# before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2 = 3;   3 = 5;   4 = 7;   5 = 6;   6 = 2;   7 = 8;   8 = 4;
    9 = 17;  10 = 10; 11 = 27; 12 = 29; 13 = 14; 14 = 9;  15 = 13;
    16 = 31; 17 = 16; 18 = 25; 19 = 20; 20 = 18; 21 = 30; 22 = 23;
    23 = 32; 24 = 11; 25 = 19; 26 = 21; 27 = 12; 28 = 22; 29 = 15;
    30 = 26; 31 = 24; 32 = 28;
    33 = 34; 34 = 35; 35 = 33;
    36 = 38; 37 = 40; 38 = 37; 39 = 42; 40 = 36; 41 = 39; 42 = 41;
    43 = 43; 44 = 45; 45 = 44;
    46 = 46; 47 = 47; 48 = 50; 49 = 48; 50 = 49;
}

$firstRow = 2
$lastRow = 50

# Snapshot columns B, C, D, E, F for every data row before mutating anything,
# so that overlapping source/destination writes don't clobber each other.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = @{
        B = $ws.Cells.Item($r, 2).Value2
        C = $ws.Cells.Item($r, 3).Value2
        D = $ws.Cells.Item($r, 4).Value2
        E = $ws.Cells.Item($r, 5).Value2
        F = $ws.Cells.Item($r, 6).Value2
    }
}

for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $mapping[$destRow]
    $src = $snapshot[$srcRow]

    $ws.Cells.Item($destRow, 2).Value2 = $src.B
    $ws.Cells.Item($destRow, 3).Value2 = $src.C
    $ws.Cells.Item($destRow, 4).Value2 = $src.D
    $ws.Cells.Item($destRow, 5).Value2 = $src.E
    $ws.Cells.Item($destRow, 6).Value2 = $src.F
}
